$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.155.25'
$ws.Range("E2").Value = '  -0.62%  '
$ws.Range("D3").Value = '2.610.80'
$ws.Range("E3").Value = '  +0.51%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '520.31'
$ws.Range("E5").Value = '  +0.92%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.75'
$ws.Range("E6").Value = '  -3.72%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("E8").Value = '  -5.12%  '
$ws.Range("D9").Value = '2.615.00'
$ws.Range("E9").Value = '  +0.23%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.30'
$ws.Range("E10").Value = '  -4.94%  '
$ws.Range("E11").Value = '  +0.22%  '
$ws.Range("E12").Value = '  -1.69%  '
$ws.Range("E13").Value = '  -0.79%  '
$ws.Range("D14").Value = '3.068.45'
$ws.Range("E14").Value = '  +0.52%  '
$ws.Range("D15").Value = '60.188.53'
$ws.Range("E15").Value = '  -0.64%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.30'
$ws.Range("E16").Value = '  -2.26%  '
$ws.Range("E17").Value = '  -1.77%  '
$ws.Range("D18").Value = '2.613.91'
$ws.Range("E18").Value = '  +0.56%  '
$ws.Range("E19").Value = '  -2.74%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '341.15'
$ws.Range("E20").Value = '  -3.48%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.40'
$ws.Range("E21").Value = '  -1.85%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.09'
$ws.Range("E22").Value = '  -2.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.992'
$ws.Range("E23").Value = '  -0.66%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.19'
$ws.Range("E24").Value = '  -1.38%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.418'
$ws.Range("E25").Value = '  -2.57%  '
$ws.Range("B26").Value = 'Binance-PegBSC-USD'
$ws.Range("C26").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.996'
$ws.Range("E26").Value = '  +7.11%  '
$ws.Range("B27").Value = 'Kaspa'
$ws.Range("C27").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.162'
$ws.Range("E27").Value = '  -2.25%  '
$ws.Range("D28").Value = '0.0₃0804'
$ws.Range("E28").Value = '  -4.32%  '
$ws.Range("E29").Value = '  -3.65%  '
$ws.Range("E30").Value = '  +0.05%  '
$ws.Range("E31").Value = '  -3.27%  '
$ws.Range("E32").Value = '  +0.12%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.91'
$ws.Range("E33").Value = '  -2.73%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '150.12'
$ws.Range("E34").Value = '  +0.03%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.95'
$ws.Range("E35").Value = '  -3.04%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.907'
$ws.Range("E36").Value = '  -3.91%  '
$ws.Range("E37").Value = '  -5.06%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.862'
$ws.Range("E38").Value = '  +2.54%  '
$ws.Range("E39").Value = '  +0.03%  '
$ws.Range("E40").Value = '  -3.42%  '
$ws.Range("E41").Value = '  -4.05%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '286.67'
$ws.Range("E43").Value = '  -0.13%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.100'
$ws.Range("E44").Value = '  -1.25%  '
$ws.Range("E45").Value = '  +0.20%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0545'
$ws.Range("E46").Value = '  -2.46%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '19.47'
$ws.Range("E47").Value = '  -0.20%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0232'
$ws.Range("E48").Value = '  -1.62%  '
$ws.Range("E49").Value = '  +0.92%  '
$ws.Range("E50").Value = '  -5.76%  '
$ws.Range("D51").Value = '1.953.19'
$ws.Range("E51").Value = '  -1.20%  '
